$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 815.6667  # ALC H129: 804.5 -> 815.6667
$ws.Cells.Item(129, 10).Value = 921.63635  # ALC J129: 944.625 -> 921.63635
$ws.Cells.Item(129, 12).Value = 2764.90905  # ALC L129: 2833.875 -> 2764.90905
$ws.Cells.Item(129, 14).Value = -12764.90905  # ALC N129: -12833.875 -> -12764.90905

$ws.Cells.Item(132, 8).Value = 6566  # ALC H132: 6921.5 -> 6566
$ws.Cells.Item(132, 9).Value = 6985  # ALC I132: 6921.5 -> 6985
$ws.Cells.Item(132, 10).Value = 700  # ALC J132: 0 -> 700
$ws.Cells.Item(132, 11).Value = 20955  # ALC K132: 20764.5 -> 20955
$ws.Cells.Item(132, 12).Value = 2100  # ALC L132: 0 -> 2100
$ws.Cells.Item(132, 13).Value = -18425  # ALC M132: -18234.5 -> -18425
$ws.Cells.Item(132, 14).Value = -7160  # ALC N132: None -> -7160

$ws.Cells.Item(135, 8).Value = 15146703  # ALC H135: 11780843 -> 15146703
$ws.Cells.Item(135, 9).Value = 4334.1816  # ALC I135: 3337.1035 -> 4334.1816
$ws.Cells.Item(135, 10).Value = 40772252  # ALC J135: 33127572 -> 40772252
$ws.Cells.Item(135, 11).Value = 39007.6344  # ALC K135: 30033.9315 -> 39007.6344
$ws.Cells.Item(135, 12).Value = 366950268  # ALC L135: 298148148 -> 366950268
$ws.Cells.Item(135, 13).Value = -36472.6344  # ALC M135: -27498.9315 -> -36472.6344
$ws.Cells.Item(135, 14).Value = -366955338  # ALC N135: -298153218 -> -366955338

$ws.Cells.Item(138, 8).Value = 3266.2346  # ALC H138: 3384.9167 -> 3266.2346
$ws.Cells.Item(138, 9).Value = 2584.5  # ALC I138: 2851.2058 -> 2584.5
$ws.Cells.Item(138, 10).Value = 3868.6978  # ALC J138: 4082.8462 -> 3868.6978
$ws.Cells.Item(138, 11).Value = 7753.5  # ALC K138: 8553.617400000001 -> 7753.5
$ws.Cells.Item(138, 12).Value = 11606.0934  # ALC L138: 12248.5386 -> 11606.0934
$ws.Cells.Item(138, 13).Value = -2613.5  # ALC M138: -3413.617400000001 -> -2613.5
$ws.Cells.Item(138, 14).Value = -21886.0934  # ALC N138: -22528.5386 -> -21886.0934

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16363.597  # ARM H32: 15558.682 -> 16363.597
$ws.Cells.Item(32, 9).Value = 13529.54  # ARM I32: 12962.863 -> 13529.54
$ws.Cells.Item(32, 10).Value = 61000  # ARM J32: 72666.664 -> 61000
$ws.Cells.Item(32, 11).Value = 13529.54  # ARM K32: 12962.863 -> 13529.54
$ws.Cells.Item(32, 12).Value = 61000  # ARM L32: 72666.664 -> 61000
$ws.Cells.Item(32, 13).Value = -13242.54  # ARM M32: -12675.863 -> -13242.54
$ws.Cells.Item(32, 14).Value = -61574  # ARM N32: -73240.664 -> -61574

$ws.Cells.Item(61, 8).Value = 3707343  # ARM H61: 3475637.8 -> 3707343
$ws.Cells.Item(61, 9).Value = 5851136.5  # ARM I61: 5293858 -> 5851136.5
$ws.Cells.Item(61, 10).Value = 4426.727  # ARM J61: 4490.364 -> 4426.727
$ws.Cells.Item(61, 11).Value = 5851136.5  # ARM K61: 5293858 -> 5851136.5
$ws.Cells.Item(61, 12).Value = 4426.727  # ARM L61: 4490.364 -> 4426.727
$ws.Cells.Item(61, 13).Value = -5850924.5  # ARM M61: -5293646 -> -5850924.5
$ws.Cells.Item(61, 14).Value = -4850.727  # ARM N61: -4914.364 -> -4850.727

$ws.Cells.Item(102, 8).Value = 3451.8333  # ARM H102: 2192.4285 -> 3451.8333
$ws.Cells.Item(102, 9).Value = 3860  # ARM I102: 2252.5386 -> 3860
$ws.Cells.Item(102, 11).Value = 3860  # ARM K102: 2252.5386 -> 3860
$ws.Cells.Item(102, 13).Value = -2238  # ARM M102: -630.5385999999999 -> -2238

$ws.Cells.Item(117, 8).Value = 20045  # ARM H117: 20060 -> 20045
$ws.Cells.Item(117, 10).Value = 20045  # ARM J117: 20060 -> 20045
$ws.Cells.Item(117, 12).Value = 20045  # ARM L117: 20060 -> 20045
$ws.Cells.Item(117, 14).Value = -29223  # ARM N117: -29238 -> -29223

$ws.Cells.Item(136, 8).Value = 3707343  # ARM H136: 3475637.8 -> 3707343
$ws.Cells.Item(136, 9).Value = 5851136.5  # ARM I136: 5293858 -> 5851136.5
$ws.Cells.Item(136, 10).Value = 4426.727  # ARM J136: 4490.364 -> 4426.727
$ws.Cells.Item(136, 11).Value = 17553409.5  # ARM K136: 15881574 -> 17553409.5
$ws.Cells.Item(136, 12).Value = 13280.181  # ARM L136: 13471.092 -> 13280.181
$ws.Cells.Item(136, 13).Value = -17550859.5  # ARM M136: -15879024 -> -17550859.5
$ws.Cells.Item(136, 14).Value = -18380.181  # ARM N136: -18571.092 -> -18380.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 2478.2222  # BSM H64: 5297 -> 2478.2222
$ws.Cells.Item(64, 9).Value = 10068  # BSM I64: 20000 -> 10068
$ws.Cells.Item(64, 10).Value = 309.7143  # BSM J64: 396 -> 309.7143
$ws.Cells.Item(64, 11).Value = 10068  # BSM K64: 20000 -> 10068
$ws.Cells.Item(64, 12).Value = 309.7143  # BSM L64: 396 -> 309.7143
$ws.Cells.Item(64, 13).Value = -9843  # BSM M64: -19775 -> -9843
$ws.Cells.Item(64, 14).Value = -759.7143  # BSM N64: -846 -> -759.7143

$ws.Cells.Item(67, 8).Value = 2478.2222  # BSM H67: 5297 -> 2478.2222
$ws.Cells.Item(67, 9).Value = 10068  # BSM I67: 20000 -> 10068
$ws.Cells.Item(67, 10).Value = 309.7143  # BSM J67: 396 -> 309.7143
$ws.Cells.Item(67, 11).Value = 10068  # BSM K67: 20000 -> 10068
$ws.Cells.Item(67, 12).Value = 309.7143  # BSM L67: 396 -> 309.7143
$ws.Cells.Item(67, 13).Value = -9288  # BSM M67: -19220 -> -9288
$ws.Cells.Item(67, 14).Value = -1869.7143  # BSM N67: -1956 -> -1869.7143

$ws.Cells.Item(86, 8).Value = 2116.6667  # BSM H86: 1795.409 -> 2116.6667
$ws.Cells.Item(86, 9).Value = 1766.6666  # BSM I86: 1683.3334 -> 1766.6666
$ws.Cells.Item(86, 10).Value = 2466.6667  # BSM J86: 2299.75 -> 2466.6667
$ws.Cells.Item(86, 11).Value = 1766.6666  # BSM K86: 1683.3334 -> 1766.6666
$ws.Cells.Item(86, 12).Value = 2466.6667  # BSM L86: 2299.75 -> 2466.6667
$ws.Cells.Item(86, 13).Value = -643.6666  # BSM M86: -560.3334 -> -643.6666
$ws.Cells.Item(86, 14).Value = -4712.6667  # BSM N86: -4545.75 -> -4712.6667

$ws.Cells.Item(89, 8).Value = 2116.6667  # BSM H89: 1795.409 -> 2116.6667
$ws.Cells.Item(89, 9).Value = 1766.6666  # BSM I89: 1683.3334 -> 1766.6666
$ws.Cells.Item(89, 10).Value = 2466.6667  # BSM J89: 2299.75 -> 2466.6667
$ws.Cells.Item(89, 11).Value = 8833.333000000001  # BSM K89: 8416.666999999999 -> 8833.333000000001
$ws.Cells.Item(89, 12).Value = 12333.3335  # BSM L89: 11498.75 -> 12333.3335
$ws.Cells.Item(89, 13).Value = -3217.333000000001  # BSM M89: -2800.666999999999 -> -3217.333000000001
$ws.Cells.Item(89, 14).Value = -23565.3335  # BSM N89: -22730.75 -> -23565.3335

$ws.Cells.Item(141, 8).Value = 56840.5  # BSM H141: 57813.547 -> 56840.5
$ws.Cells.Item(141, 10).Value = 52737.7  # BSM J141: 53471.11 -> 52737.7
$ws.Cells.Item(141, 12).Value = 52737.7  # BSM L141: 53471.11 -> 52737.7
$ws.Cells.Item(141, 14).Value = -63097.7  # BSM N141: -63831.11 -> -63097.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 205.6  # CRP H22: 204.11111 -> 205.6
$ws.Cells.Item(22, 9).Value = 199.56522  # CRP I22: 196.44 -> 199.56522
$ws.Cells.Item(22, 10).Value = 275  # CRP J22: 300 -> 275
$ws.Cells.Item(22, 11).Value = 199.56522  # CRP K22: 196.44 -> 199.56522
$ws.Cells.Item(22, 12).Value = 275  # CRP L22: 300 -> 275
$ws.Cells.Item(22, 13).Value = 150.43478  # CRP M22: 153.56 -> 150.43478
$ws.Cells.Item(22, 14).Value = -975  # CRP N22: -1000 -> -975

$ws.Cells.Item(93, 8).Value = 9659.5  # CRP H93: 9751.166999999999 -> 9659.5
$ws.Cells.Item(93, 9).Value = 7810.364  # CRP I93: 7910.364 -> 7810.364
$ws.Cells.Item(93, 11).Value = 7810.364  # CRP K93: 7910.364 -> 7810.364
$ws.Cells.Item(93, 13).Value = -5938.364  # CRP M93: -6038.364 -> -5938.364

$ws.Cells.Item(134, 8).Value = 2383  # CRP H134: 1986.4 -> 2383
$ws.Cells.Item(134, 9).Value = 2142.739  # CRP I134: 1795.6897 -> 2142.739
$ws.Cells.Item(134, 10).Value = 3073.75  # CRP J134: 2489.182 -> 3073.75
$ws.Cells.Item(134, 11).Value = 6428.217000000001  # CRP K134: 5387.0691 -> 6428.217000000001
$ws.Cells.Item(134, 12).Value = 9221.25  # CRP L134: 7467.545999999999 -> 9221.25
$ws.Cells.Item(134, 13).Value = -3893.217000000001  # CRP M134: -2852.0691 -> -3893.217000000001
$ws.Cells.Item(134, 14).Value = -14291.25  # CRP N134: -12537.546 -> -14291.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 150.66667  # CUL H40: 270.66666 -> 150.66667
$ws.Cells.Item(40, 9).Value = 150.66667  # CUL I40: 182 -> 150.66667
$ws.Cells.Item(40, 10).Value = 0  # CUL J40: 980 -> 0
$ws.Cells.Item(40, 11).Value = 602.66668  # CUL K40: 728 -> 602.66668
$ws.Cells.Item(40, 12).Value = 0  # CUL L40: 3920 -> 0
$ws.Cells.Item(40, 13).Value = -533.66668  # CUL M40: -659 -> -533.66668
$ws.Cells.Item(40, 14).ClearContents()  # CUL N40: -4058 -> (removed)

$ws.Cells.Item(56, 8).Value = 3534.8  # CUL H56: 3562.9473 -> 3534.8
$ws.Cells.Item(56, 9).Value = 3534.8  # CUL I56: 3562.9473 -> 3534.8
$ws.Cells.Item(56, 11).Value = 3534.8  # CUL K56: 3562.9473 -> 3534.8
$ws.Cells.Item(56, 13).Value = -3004.8  # CUL M56: -3032.9473 -> -3004.8

$ws.Cells.Item(86, 8).Value = 756.7143  # CUL H86: 771 -> 756.7143
$ws.Cells.Item(86, 10).Value = 899.8  # CUL J86: 919.8 -> 899.8
$ws.Cells.Item(86, 12).Value = 2699.4  # CUL L86: 2759.4 -> 2699.4
$ws.Cells.Item(86, 14).Value = -5071.4  # CUL N86: -5131.4 -> -5071.4

$ws.Cells.Item(89, 8).Value = 756.7143  # CUL H89: 771 -> 756.7143
$ws.Cells.Item(89, 10).Value = 899.8  # CUL J89: 919.8 -> 899.8
$ws.Cells.Item(89, 12).Value = 8098.2  # CUL L89: 8278.199999999999 -> 8098.2
$ws.Cells.Item(89, 14).Value = -19954.2  # CUL N89: -20134.2 -> -19954.2

$ws.Cells.Item(113, 8).Value = 594.3333  # CUL H113: 595.9048 -> 594.3333
$ws.Cells.Item(113, 10).Value = 564.1579  # CUL J113: 567.6316 -> 564.1579
$ws.Cells.Item(113, 12).Value = 1692.4737  # CUL L113: 1702.8948 -> 1692.4737
$ws.Cells.Item(113, 14).Value = -6032.4737  # CUL N113: -6042.8948 -> -6032.4737

$ws.Cells.Item(132, 8).Value = 52632696  # CUL H132: 58824750 -> 52632696
$ws.Cells.Item(132, 9).Value = 76924024  # CUL I132: 100001020 -> 76924024
$ws.Cells.Item(132, 10).Value = 1479.6666  # CUL J132: 1511.1428 -> 1479.6666
$ws.Cells.Item(132, 11).Value = 692316216  # CUL K132: 900009180 -> 692316216
$ws.Cells.Item(132, 12).Value = 13316.9994  # CUL L132: 13600.2852 -> 13316.9994
$ws.Cells.Item(132, 13).Value = -692313686  # CUL M132: -900006650 -> -692313686
$ws.Cells.Item(132, 14).Value = -18376.9994  # CUL N132: -18660.2852 -> -18376.9994

$ws.Cells.Item(137, 8).Value = 25423.367  # CUL H137: 26718.412 -> 25423.367
$ws.Cells.Item(137, 9).Value = 5490  # CUL I137: 2952 -> 5490
$ws.Cells.Item(137, 10).Value = 27195.223  # CUL J137: 32499.432 -> 27195.223
$ws.Cells.Item(137, 11).Value = 16470  # CUL K137: 8856 -> 16470
$ws.Cells.Item(137, 12).Value = 81585.66900000001  # CUL L137: 97498.296 -> 81585.66900000001
$ws.Cells.Item(137, 13).Value = -11370  # CUL M137: -3756 -> -11370
$ws.Cells.Item(137, 14).Value = -91785.66900000001  # CUL N137: -107698.296 -> -91785.66900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4662.857  # GSM H70: 4704.303 -> 4662.857
$ws.Cells.Item(70, 9).Value = 4482.8335  # GSM I70: 4569.5884 -> 4482.8335
$ws.Cells.Item(70, 10).Value = 4853.4707  # GSM J70: 4847.4375 -> 4853.4707
$ws.Cells.Item(70, 11).Value = 4482.8335  # GSM K70: 4569.5884 -> 4482.8335
$ws.Cells.Item(70, 12).Value = 4853.4707  # GSM L70: 4847.4375 -> 4853.4707
$ws.Cells.Item(70, 13).Value = -4212.8335  # GSM M70: -4299.5884 -> -4212.8335
$ws.Cells.Item(70, 14).Value = -5393.4707  # GSM N70: -5387.4375 -> -5393.4707

$ws.Cells.Item(73, 8).Value = 4662.857  # GSM H73: 4704.303 -> 4662.857
$ws.Cells.Item(73, 9).Value = 4482.8335  # GSM I73: 4569.5884 -> 4482.8335
$ws.Cells.Item(73, 10).Value = 4853.4707  # GSM J73: 4847.4375 -> 4853.4707
$ws.Cells.Item(73, 11).Value = 4482.8335  # GSM K73: 4569.5884 -> 4482.8335
$ws.Cells.Item(73, 12).Value = 4853.4707  # GSM L73: 4847.4375 -> 4853.4707
$ws.Cells.Item(73, 13).Value = -3546.8335  # GSM M73: -3633.5884 -> -3546.8335
$ws.Cells.Item(73, 14).Value = -6725.4707  # GSM N73: -6719.4375 -> -6725.4707

$ws.Cells.Item(80, 8).Value = 7897.1  # GSM H80: 7591.4287 -> 7897.1
$ws.Cells.Item(80, 9).Value = 2542.5  # GSM I80: 2482.2222 -> 2542.5
$ws.Cells.Item(80, 10).Value = 11466.833  # GSM J80: 11423.333 -> 11466.833
$ws.Cells.Item(80, 11).Value = 2542.5  # GSM K80: 2482.2222 -> 2542.5
$ws.Cells.Item(80, 12).Value = 11466.833  # GSM L80: 11423.333 -> 11466.833
$ws.Cells.Item(80, 13).Value = -1544.5  # GSM M80: -1484.2222 -> -1544.5
$ws.Cells.Item(80, 14).Value = -13462.833  # GSM N80: -13419.333 -> -13462.833

$ws.Cells.Item(83, 8).Value = 7897.1  # GSM H83: 7591.4287 -> 7897.1
$ws.Cells.Item(83, 9).Value = 2542.5  # GSM I83: 2482.2222 -> 2542.5
$ws.Cells.Item(83, 10).Value = 11466.833  # GSM J83: 11423.333 -> 11466.833
$ws.Cells.Item(83, 11).Value = 12712.5  # GSM K83: 12411.111 -> 12712.5
$ws.Cells.Item(83, 12).Value = 57334.165  # GSM L83: 57116.665 -> 57334.165
$ws.Cells.Item(83, 13).Value = -7720.5  # GSM M83: -7419.111000000001 -> -7720.5
$ws.Cells.Item(83, 14).Value = -67318.16500000001  # GSM N83: -67100.66500000001 -> -67318.16500000001

$ws.Cells.Item(97, 8).Value = 2133.3333  # GSM H97: 2066.6667 -> 2133.3333
$ws.Cells.Item(97, 9).Value = 1920  # GSM I97: 1928.5714 -> 1920
$ws.Cells.Item(97, 10).Value = 3200  # GSM J97: 2550 -> 3200
$ws.Cells.Item(97, 11).Value = 1920  # GSM K97: 1928.5714 -> 1920
$ws.Cells.Item(97, 12).Value = 3200  # GSM L97: 2550 -> 3200
$ws.Cells.Item(97, 13).Value = -1424  # GSM M97: -1432.5714 -> -1424
$ws.Cells.Item(97, 14).Value = -4192  # GSM N97: -3542 -> -4192

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 7019.077  # LTW H132: 8509.700000000001 -> 7019.077
$ws.Cells.Item(132, 9).Value = 7690.381  # LTW I132: 8788.666999999999 -> 7690.381
$ws.Cells.Item(132, 10).Value = 4199.6  # LTW J132: 5999 -> 4199.6
$ws.Cells.Item(132, 11).Value = 23071.143  # LTW K132: 26366.001 -> 23071.143
$ws.Cells.Item(132, 12).Value = 12598.8  # LTW L132: 17997 -> 12598.8
$ws.Cells.Item(132, 13).Value = -20541.143  # LTW M132: -23836.001 -> -20541.143
$ws.Cells.Item(132, 14).Value = -17658.8  # LTW N132: -23057 -> -17658.8

$ws.Cells.Item(136, 8).Value = 2125.625  # LTW H136: 2207.4 -> 2125.625
$ws.Cells.Item(136, 9).Value = 1167.5  # LTW I136: 1175.9166 -> 1167.5
$ws.Cells.Item(136, 10).Value = 5000  # LTW J136: 6333.3335 -> 5000
$ws.Cells.Item(136, 11).Value = 3502.5  # LTW K136: 3527.7498 -> 3502.5
$ws.Cells.Item(136, 12).Value = 15000  # LTW L136: 19000.0005 -> 15000
$ws.Cells.Item(136, 13).Value = -952.5  # LTW M136: -977.7498000000001 -> -952.5
$ws.Cells.Item(136, 14).Value = -20100  # LTW N136: -24100.0005 -> -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3231.842  # WVR H132: 3267.838 -> 3231.842
$ws.Cells.Item(132, 9).Value = 3144.8965  # WVR I132: 3165.5862 -> 3144.8965
$ws.Cells.Item(132, 10).Value = 3512  # WVR J132: 3638.5 -> 3512
$ws.Cells.Item(132, 11).Value = 9434.6895  # WVR K132: 9496.758600000001 -> 9434.6895
$ws.Cells.Item(132, 12).Value = 10536  # WVR L132: 10915.5 -> 10536
$ws.Cells.Item(132, 13).Value = -6904.6895  # WVR M132: -6966.758600000001 -> -6904.6895
$ws.Cells.Item(132, 14).Value = -15596  # WVR N132: -15975.5 -> -15596

$ws.Cells.Item(136, 8).Value = 7954.3784  # WVR H136: 7835.65 -> 7954.3784
$ws.Cells.Item(136, 9).Value = 12007.363  # WVR I136: 10170.214 -> 12007.363
$ws.Cells.Item(136, 10).Value = 2010  # WVR J136: 2388.3333 -> 2010
$ws.Cells.Item(136, 11).Value = 36022.089  # WVR K136: 30510.642 -> 36022.089
$ws.Cells.Item(136, 12).Value = 6030  # WVR L136: 7164.999899999999 -> 6030
$ws.Cells.Item(136, 13).Value = -33472.089  # WVR M136: -27960.642 -> -33472.089
$ws.Cells.Item(136, 14).Value = -11130  # WVR N136: -12264.9999 -> -11130
